$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirement")

# --- Add the new "DegreeId" column (D) -------------------------------------
$ws.Range("D1").Value = "DegreeId"
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
}

# --- Update the formulas in column F to also emit the DegreeId ------------
# F2 carries its own (non-shared) formula.
$ws.Range("F2").Formula = '=$E$2&$A$1&"="&A2&","&$B$1&"="&"''"&B2&"''"&","&$C$1&"="&"''"&C2&"''"&","&$D$1&"="&D2&$G$2'

# F3:F14 are a shared-formula block anchored at F3.
$ws.Range("F3:F14").Formula = '=$E$2&$A$1&"="&A3&","&$B$1&"="&"''"&B3&"''"&","&$C$1&"="&"''"&C3&"''"&","&$D$1&"="&D3&$G$2'

# --- Column widths: new column E (bestfit) and wider column F -------------
$ws.Columns.Item(5).ColumnWidth = 17.166666666666668
$ws.Columns.Item(6).ColumnWidth = 123.66666666666667

# --- View / selection state --------------------------------------------
# The "Requirement" sheet becomes the active tab/sheet, scrolled right a
# little, with F2:F14 selected.
$ws.Activate() | Out-Null
$ws.Range("F2:F14").Select() | Out-Null
